$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, pushing the existing rows 72-76 down to 73-77.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly price entry.
$ws.Range("A72").Value = 1
$ws.Range("B72").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C72").Value = "Arica y Parinacota"
$ws.Range("D72").Value = 44858
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = 100112009
$ws.Range("G72").Value = "Acelga"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 200
$ws.Range("K72").Value = 2800
$ws.Range("L72").Value = 3000
$ws.Range("M72").Value = 2900
$ws.Range("N72").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O72").Value = "Región de Arica y Parinacota"
$ws.Range("P72").Value = 967
$ws.Range("Q72").Value = 3
$ws.Range("R72").Value = "Hortaliza"
